$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.839.43"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.314.16"
$ws.Range("E3").Value = "  +3.15%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "97.59"
$ws.Range("E5").Value = "  +2.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "272.00"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.625"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.18"
$ws.Range("E10").Value = "  -2.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0952"
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.02"
$ws.Range("E12").Value = "  -3.19%  "
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.652.55"
$ws.Range("E14").Value = "  +2.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.46"
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.872"
$ws.Range("E16").Value = "  +6.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.320.42"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.796.09"
$ws.Range("E18").Value = "  -0.32%  "
$ws.Range("E19").Value = "  +3.93%  "
$ws.Range("E20").Value = "  +4.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.49"
$ws.Range("E21").Value = "  +3.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "239.72"
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.43"
$ws.Range("E24").Value = "  +2.57%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.36"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.49"
$ws.Range("E28").Value = "  -2.73%  "
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.04"
$ws.Range("E30").Value = "  -5.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.37"
$ws.Range("E31").Value = "  +6.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "175.50"
$ws.Range("E32").Value = "  +1.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0911"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.48"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("E36").Value = "  +3.01%  "
$ws.Range("E37").Value = "  -3.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.44"
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("E39").Value = "  -7.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.244"
$ws.Range("E40").Value = "  +7.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.40"
$ws.Range("E41").Value = "  +10.06%  "
$ws.Range("E42").Value = "  +22.13%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.34"
$ws.Range("E43").Value = "  -4.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.80"
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.09"
$ws.Range("E45").Value = "  +8.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.33"
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("E47").Value = "  +3.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.34"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("E50").Value = "  +15.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.539.84"
$ws.Range("E51").Value = "  +3.03%  "
